$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the header row (columns stay in A:E, but the labels/order change)
$ws.Range("A1").Value = "WORD"
$ws.Range("B1").Value = "TRANSLATION"
$ws.Range("C1").Value = "CATEGORY"
$ws.Range("D1").Value = "LANGUAGE"
$ws.Range("E1").Value = "LEVEL"

# Existing placeholder row becomes real-ish sample data
$ws.Range("A2").Value = "aaa"
$ws.Range("B2").Value = "aaa"
$ws.Range("C2").Value = "aaa"
$ws.Range("D2").Value = "ES"
$ws.Range("E2").Value = 5

# New row 3
$ws.Range("A3").Value = "agotado"
$ws.Range("B3").Value = "ausverkauft/vergriffen"
$ws.Range("C3").Value = "adjective"
$ws.Range("D3").Value = "ES"
$ws.Range("E3").Value = 0

# New row 4
$ws.Range("A4").Value = "alegre"
$ws.Range("B4").Value = "fröhlich"
$ws.Range("C4").Value = "adjective"
$ws.Range("D4").Value = "ES"
$ws.Range("E4").Value = 0
